$wb = $excel.ActiveWorkbook

# This script applies updated profit-calculation figures (recomputed by the
# scheduled market-data runner) to the per-sheet Leve profit tables.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 16928.143
$ws.Range("I21").Value = 5000
$ws.Range("J21").Value = 18916.166
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 18916.166
$ws.Range("M21").Value = -4532
$ws.Range("N21").Value = -19852.166

$ws.Range("H23").Value = 16928.143
$ws.Range("I23").Value = 5000
$ws.Range("J23").Value = 18916.166
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 18916.166
$ws.Range("M23").Value = -4766
$ws.Range("N23").Value = -19384.166

$ws.Range("H33").Value = 312.53333
$ws.Range("I33").Value = 330.39285
$ws.Range("J33").Value = 62.5
$ws.Range("K33").Value = 330.39285
$ws.Range("L33").Value = 62.5
$ws.Range("M33").Value = -101.39285
$ws.Range("N33").Value = -520.5

$ws.Range("H34").Value = 896
$ws.Range("I34").Value = 896
$ws.Range("K34").Value = 896
$ws.Range("M34").Value = -693

$ws.Range("H36").Value = 896
$ws.Range("I36").Value = 896
$ws.Range("K36").Value = 896
$ws.Range("M36").Value = -181

$ws.Range("H38").Value = 451.46667
$ws.Range("I38").Value = 263.55554
$ws.Range("J38").Value = 733.3333
$ws.Range("K38").Value = 790.66662
$ws.Range("L38").Value = 2199.9999
$ws.Range("M38").Value = -418.66662
$ws.Range("N38").Value = -2943.9999

$ws.Range("H58").Value = 1351.3636
$ws.Range("I58").Value = 1207.2222
$ws.Range("K58").Value = 3621.6666
$ws.Range("M58").Value = -3471.6666

$ws.Range("H87").Value = 24999.62
$ws.Range("J87").Value = 24999.62
$ws.Range("L87").Value = 24999.62
$ws.Range("N87").Value = -27495.62

$ws.Range("H90").Value = 24999.62
$ws.Range("J90").Value = 24999.62
$ws.Range("L90").Value = 74998.86
$ws.Range("N90").Value = -87478.86

$ws.Range("H107").Value = 1015.5
$ws.Range("I107").Value = 1015.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1015.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 904.5
$ws.Range("N107").ClearContents()

$ws.Range("H129").Value = 1344.5428
$ws.Range("I129").Value = 607.3333
$ws.Range("J129").Value = 1729.174
$ws.Range("K129").Value = 1821.9999
$ws.Range("L129").Value = 5187.522
$ws.Range("M129").Value = 3178.0001
$ws.Range("N129").Value = -15187.522

$ws.Range("H137").Value = 1363.2034
$ws.Range("I137").Value = 1275.2903
$ws.Range("J137").Value = 1460.5358
$ws.Range("K137").Value = 3825.8709
$ws.Range("L137").Value = 4381.607400000001
$ws.Range("M137").Value = -1275.8709
$ws.Range("N137").Value = -9481.607400000001

$ws.Range("H138").Value = 1564.55
$ws.Range("I138").Value = 686.5897
$ws.Range("J138").Value = 2125.869
$ws.Range("K138").Value = 2059.7691
$ws.Range("L138").Value = 6377.607
$ws.Range("M138").Value = 3080.2309
$ws.Range("N138").Value = -16657.607


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1482.0476
$ws.Range("I2").Value = 1512.3529
$ws.Range("J2").Value = 1353.25
$ws.Range("K2").Value = 1512.3529
$ws.Range("L2").Value = 1353.25
$ws.Range("M2").Value = -1399.3529
$ws.Range("N2").Value = -1579.25

$ws.Range("H32").Value = 7108.41
$ws.Range("I32").Value = 5105
$ws.Range("J32").Value = 21800.084
$ws.Range("K32").Value = 5105
$ws.Range("L32").Value = 21800.084
$ws.Range("M32").Value = -4818
$ws.Range("N32").Value = -22374.084

$ws.Range("H61").Value = 6412254.5
$ws.Range("I61").Value = 7753603
$ws.Range("J61").Value = 3588.889
$ws.Range("K61").Value = 7753603
$ws.Range("L61").Value = 3588.889
$ws.Range("M61").Value = -7753391
$ws.Range("N61").Value = -4012.889

$ws.Range("H110").Value = 113562
$ws.Range("I110").Value = 113562
$ws.Range("K110").Value = 113562
$ws.Range("M110").Value = -111517

$ws.Range("H116").Value = 1482.0476
$ws.Range("I116").Value = 1512.3529
$ws.Range("J116").Value = 1353.25
$ws.Range("K116").Value = 1512.3529
$ws.Range("L116").Value = 1353.25
$ws.Range("M116").Value = 781.6470999999999
$ws.Range("N116").Value = -5941.25

$ws.Range("H136").Value = 6412254.5
$ws.Range("I136").Value = 7753603
$ws.Range("J136").Value = 3588.889
$ws.Range("K136").Value = 23260809
$ws.Range("L136").Value = 10766.667
$ws.Range("M136").Value = -23258259
$ws.Range("N136").Value = -15866.667


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1482.0476
$ws.Range("I3").Value = 1512.3529
$ws.Range("J3").Value = 1353.25
$ws.Range("K3").Value = 1512.3529
$ws.Range("L3").Value = 1353.25
$ws.Range("M3").Value = -1398.3529
$ws.Range("N3").Value = -1581.25


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4233.087
$ws.Range("I31").Value = 1338.6578
$ws.Range("J31").Value = 7781.0967
$ws.Range("K31").Value = 1338.6578
$ws.Range("L31").Value = 7781.0967
$ws.Range("M31").Value = -1043.6578
$ws.Range("N31").Value = -8371.0967

$ws.Range("H34").Value = 4233.087
$ws.Range("I34").Value = 1338.6578
$ws.Range("J34").Value = 7781.0967
$ws.Range("K34").Value = 1338.6578
$ws.Range("L34").Value = 7781.0967
$ws.Range("M34").Value = -1136.6578
$ws.Range("N34").Value = -8185.0967

$ws.Range("H58").Value = 1112.8422
$ws.Range("I58").Value = 863.6316
$ws.Range("J58").Value = 1611.2632
$ws.Range("K58").Value = 863.6316
$ws.Range("L58").Value = 1611.2632
$ws.Range("M58").Value = -660.6316
$ws.Range("N58").Value = -2017.2632

$ws.Range("H136").Value = 1112.8422
$ws.Range("I136").Value = 863.6316
$ws.Range("J136").Value = 1611.2632
$ws.Range("K136").Value = 2590.8948
$ws.Range("L136").Value = 4833.7896
$ws.Range("M136").Value = -40.89480000000003
$ws.Range("N136").Value = -9933.7896


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 13158367
$ws.Range("I34").Value = 258.75
$ws.Range("J34").Value = 14706379
$ws.Range("K34").Value = 776.25
$ws.Range("L34").Value = 44119137
$ws.Range("M34").Value = -692.25
$ws.Range("N34").Value = -44119305

$ws.Range("H39").Value = 2449.625
$ws.Range("I39").Value = 1399
$ws.Range("J39").Value = 2599.7144
$ws.Range("K39").Value = 4197
$ws.Range("L39").Value = 7799.1432
$ws.Range("M39").Value = -3903
$ws.Range("N39").Value = -8387.143199999999

$ws.Range("H55").Value = 1541.0588
$ws.Range("I55").Value = 800
$ws.Range("J55").Value = 1587.375
$ws.Range("K55").Value = 2400
$ws.Range("L55").Value = 4762.125
$ws.Range("M55").Value = -2223
$ws.Range("N55").Value = -5116.125


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2799.5356
$ws.Range("I61").Value = 2606.1428
$ws.Range("J61").Value = 3379.7144
$ws.Range("K61").Value = 2606.1428
$ws.Range("L61").Value = 3379.7144
$ws.Range("M61").Value = -2404.1428
$ws.Range("N61").Value = -3783.7144

$ws.Range("H113").Value = 2799.5356
$ws.Range("I113").Value = 2606.1428
$ws.Range("J113").Value = 3379.7144
$ws.Range("K113").Value = 2606.1428
$ws.Range("L113").Value = 3379.7144
$ws.Range("M113").Value = -436.1428000000001
$ws.Range("N113").Value = -7719.7144

$ws.Range("H132").Value = 2896.5283
$ws.Range("I132").Value = 2674.7222
$ws.Range("J132").Value = 3366.2354
$ws.Range("K132").Value = 8024.1666
$ws.Range("L132").Value = 10098.7062
$ws.Range("M132").Value = -5494.1666
$ws.Range("N132").Value = -15158.7062

$ws.Range("H136").Value = 4763470
$ws.Range("I136").Value = 1404.1482
$ws.Range("J136").Value = 20835444
$ws.Range("K136").Value = 4212.444600000001
$ws.Range("L136").Value = 62506332
$ws.Range("M136").Value = -1662.444600000001
$ws.Range("N136").Value = -62511432


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 10496.5
$ws.Range("J54").Value = 10496.5
$ws.Range("L54").Value = 10496.5
$ws.Range("N54").Value = -11536.5

$ws.Range("H81").Value = 3058.6667
$ws.Range("I81").Value = 3528.7334
$ws.Range("J81").Value = 2275.2222
$ws.Range("K81").Value = 7057.4668
$ws.Range("L81").Value = 4550.4444
$ws.Range("M81").Value = -5996.4668
$ws.Range("N81").Value = -6672.4444

$ws.Range("H84").Value = 3058.6667
$ws.Range("I84").Value = 3528.7334
$ws.Range("J84").Value = 2275.2222
$ws.Range("K84").Value = 35287.334
$ws.Range("L84").Value = 22752.222
$ws.Range("M84").Value = -29983.334
$ws.Range("N84").Value = -33360.222

$ws.Range("H113").Value = 1140.7778
$ws.Range("I113").Value = 1174.3914
$ws.Range("K113").Value = 3523.1742
$ws.Range("M113").Value = -1353.1742

$ws.Range("H132").Value = 4862794
$ws.Range("I132").Value = 1830.3243
$ws.Range("J132").Value = 12682605
$ws.Range("K132").Value = 5490.9729
$ws.Range("L132").Value = 38047815
$ws.Range("M132").Value = -2960.9729
$ws.Range("N132").Value = -38052875

$ws.Range("H136").Value = 11630433
$ws.Range("I136").Value = 8199314.5
$ws.Range("J136").Value = 20002362
$ws.Range("K136").Value = 24597943.5
$ws.Range("L136").Value = 60007086
$ws.Range("M136").Value = -24595393.5
$ws.Range("N136").Value = -60012186

